# Add report for covered pairs
# Collapse the ratings table from 14 data rows (rows 2-15) down to 6 data
# rows (rows 2-7), with updated values, so the sheet now reports the
# covered pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7 (columns B..F). Column A (the 0-based index) is
# left untouched.
$data = @(
    @(5, "ae868375-47d2-4b8a-995b-66e8e33ebf0a", "BGOjoywll3sSUYo2pBwY", 5, "2nRJa9SKRgUZBHy9Ti9w"),
    @(4, "ae868375-47d2-4b8a-995b-66e8e33ebf0a", "5dwojHXzFd2RzJ2ZPNtD", 3, "PhbCHezTkzTl3gkwCF5U"),
    @(4, "ae868375-47d2-4b8a-995b-66e8e33ebf0a", "5dwojHXzFd2RzJ2ZPNtD", 4, "YShT4GJbXiVYmuGPYfv4"),
    @(3, "ae868375-47d2-4b8a-995b-66e8e33ebf0a", "PqpflZDmB5tBiM15v9MQ", 4, "hJMVI1CTbDSS4108H6Qp"),
    @(5, "ae868375-47d2-4b8a-995b-66e8e33ebf0a", "BGOjoywll3sSUYo2pBwY", 3, "hsCX7T7tqPKrlJp6WCcH"),
    @(4, "ae868375-47d2-4b8a-995b-66e8e33ebf0a", "OyhrkDZsWy64SkqfZAY5", 2, "yHra9FHoHYK939FWtJEc")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
}

# Remove the now-unused rows 8-15 (old data rows 7-13 which no longer exist
# in the trimmed report).
$ws.Rows("8:15").Delete()
